# vocabulary.xlsx - add "observation" table rows (observation_type_concept_id,
# obs_event_field_concept_id) ahead of the existing FS / unit rows, fix a couple
# of typos/wording tweaks, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows right before the current row 69 (observation / FS block) ---
$ws.Rows.Item(69).Resize(2).Insert()

# New row 69: observation_type_concept_id -> EHRencounter
$ws.Range("A69").Value = "observation"
$ws.Range("B69").Value = "observation_type_concept_id"
$ws.Range("C69").Value = "EHRencounter"
$ws.Range("D69").Value = 32827
$ws.Range("E69").Value = "EHR encounter record"
$ws.Rows.Item(69).RowHeight = 12.8

# New row 70: obs_event_field_concept_id -> procedure_occurrence
$ws.Range("A70").Value = "observation"
$ws.Range("B70").Value = "obs_event_field_concept_id"
$ws.Range("C70").Value = "procedure_occurrence"
$ws.Range("D70").Value = 1147301
$ws.Range("E70").Value = "procedure_occurrence table"
$ws.Rows.Item(70).RowHeight = 12.8

# Rows 71-73 are the old rows 69-71, shifted down by the insert above;
# row 71 keeps its taller "wrap" height automatically (ht 23.85).

# Fix the typo in the observation_concept_id field name (was observation_concept_it)
$ws.Range("B72").Value = "observation_concept_id"

# Fix "hertz" -> "Hz" unit source term
$ws.Range("C73").Value = "Hz"

# Wording tweaks elsewhere on the sheet
$ws.Range("E56").Value = "Root Mean Squared Successive Differences (RMSSD)"
$ws.Range("E61").Value = "Ratio of HRV Low and High Frequency powers"

# Move the selection / view to match the edited area
[void]$ws.Range("E62").Select()
